$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the camera-setup rows and the input-mapping row as complete in column A.
$ws.Range("A2").Value = "已完成"
$ws.Range("A3").Value = "已完成"
$ws.Range("A4").Value = "已完成基本的行走和四处看"

# Move the active cell selection to A2, matching the author's saved view.
$ws.Range("A2").Select()
